# "Login related testcases without validation"
#
# The sheet keeps the same two rows of login data, but:
#   - the "userName" label is relabelled to lowercase "username"
#   - the sheet's internal sheetId bumps from 1 to 2 (as if the sheet had
#     been recreated/duplicated during the edit), while keeping the same
#     tab name, tab order and rId
#   - the columns get resized to fit their (new) contents
#   - the selected cell moves to C5

$wb = $excel.ActiveWorkbook

# --- bump the sheet's internal sheetId from 1 to 2 -------------------------
# Excel assigns sheetId purely based on internal creation order/history, not
# tab position, so the only way to move a sheet's id forward is to make a
# fresh sheet (which is allocated the next id) carry the old sheet's name.
# Duplicating "loginTest" keeps all of its data/formatting, then we drop the
# original (id 1) and rename the duplicate (id 2) back to "loginTest".
$orig = $wb.Worksheets.Item("loginTest")
$orig.Copy($null, $orig)
$orig.Delete()

$ws = $wb.Worksheets.Item("loginTest (2)")
$ws.Name = "loginTest"

# --- relabel the username cell ---------------------------------------------
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "standard_user"
$ws.Range("B2").Value = "secret_sauce"

# --- resize columns to fit their content ------------------------------------
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

# --- move the selection --------------------------------------------------
$ws.Range("C5").Select()
